$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column at A, shifting existing data (old A:E) to (new B:F)
$ws.Columns.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "Metodo"
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Method names in new column A
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Re-fit the (now narrower/wider) first three columns to their content, like
# Excel's own best-fit column sizing would after the insert.
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 3.1666667
$ws.Columns.Item(3).ColumnWidth = 2.3333333

$wb.Save()
